$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.5552411675453186
$ws.Range("B1").Value = 1.849376797676086
$ws.Range("C1").Value = 3.552275657653809
$ws.Range("D1").Value = 1.52636456489563
$ws.Range("E1").Value = 0.7340560555458069
